$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataset")

# Append two new data rows for the tensorflow/ranking repository
# (ids 172 and 173), mirroring the existing rows in the dataset sheet.
$newRows = @(
    @(172, "https://github.com/tensorflow/ranking", "ranking", "tensorflow", "12/03/2018", "0", "0", "0", "1", "0", "1", "0", "0", "1", "0"),
    @(173, "https://github.com/tensorflow/ranking", "ranking", "tensorflow", "12/03/2018", "0", "0", "0", "1", "1", "1", "0", "0", "1", "0")
)

$startRow = 173

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    # Column A holds the numeric row id, formatted like the rest of the
    # "id" column (bold, bordered, centered).
    $idCell = $ws.Cells.Item($r, 1)
    $idCell.Value = $rowData[0]
    $idCell.Font.Bold = $true
    $idCell.Borders.LineStyle = 1
    $idCell.HorizontalAlignment = -4108
    $idCell.VerticalAlignment = -4160

    # Columns E (date) and F-O (the "0"/"1" community-smell flags) must
    # stay plain text, matching every other row in this sheet.
    for ($c = 2; $c -le $rowData.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($c -eq 5 -or $c -ge 6) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $rowData[$c - 1]
    }
}
